# Add data for 2022-05-06:
# extends the "through April 27" rolling window to "through April 28"
# for the current year column (B), and backfills one more day's worth
# of historical carjacking counts for the matching April-28 date across
# the other tracked years' columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-04-28"

# 2. Update the header label in B1 to match (shared string text).
$ws.Range("B1").Value = "April 2022 (through April 28)"

# 3. Row 2 - Austin
$ws.Range("B2").Value = 9
$ws.Range("J2").Value = 11
$ws.Range("V2").Value = 7

# 4. Row 4 - North Lawndale
$ws.Range("N4").Value = 4

# 5. Row 13 - Wicker Park (new cell)
$ws.Range("J13").Value = 1

# 6. Row 16 - Washington Heights (new cell)
$ws.Range("V16").Value = 1

# 7. Row 25 - South Shore
$ws.Range("B25").Value = 4
$ws.Range("F25").Value = 5
$ws.Range("N25").Value = 4

# 8. Row 28 - Chatham
$ws.Range("J28").Value = 4

# 9. Row 29 - West Town
$ws.Range("B29").Value = 1
$ws.Range("V29").Value = 2

# 10. Row 43 - Bridgeport (new cell)
$ws.Range("F43").Value = 1

# 11. Row 45 - Douglas (new cell)
$ws.Range("V45").Value = 1

# 12. Row 50 - Grand Crossing
$ws.Range("F50").Value = 2

# 13. Row 54 - Logan Square
$ws.Range("B54").Value = 1
$ws.Range("V54").Value = 3
